$wb = $excel.ActiveWorkbook

# --- AIC sheet: update the AICc model-selection table (rows 18-22) ---
$ws1 = $wb.Worksheets.Item("AIC")

# Row 18 ("Era + Size" model, K=5)
$ws1.Range("B18").Value = "Era + Size"
$ws1.Range("C18").Value = 5
$ws1.Range("D18").Value = -226.780339658477
$ws1.Range("E18").Value = 0
$ws1.Range("F18").Value = 1
$ws1.Range("G18").Value = 0.37616564160315302
$ws1.Range("H18").Value = 118.598503162572
$ws1.Range("I18").Value = 0.37616564160315302

# Row 19 ("Size" model, K=4)
$ws1.Range("B19").Value = "Size"
$ws1.Range("C19").Value = 4
$ws1.Range("D19").Value = -226.449966192132
$ws1.Range("E19").Value = 0.33037346634500903
$ws1.Range("F19").Value = 0.84773538898840595
$ws1.Range("G19").Value = 0.31888892650852202
$ws1.Range("H19").Value = 117.362914130549
$ws1.Range("I19").Value = 0.69505456811167499

# Row 20 ("Era x Size" model, K=6) - values only
$ws1.Range("D20").Value = -226.36054619087301
$ws1.Range("E20").Value = 0.41979346760339797
$ws1.Range("F20").Value = 0.810667956245818
$ws1.Range("G20").Value = 0.30494543188832501
$ws1.Range("H20").Value = 119.473979389143

# Row 21 ("Era" model, K=4) - values only
$ws1.Range("D21").Value = -13.537492937610301
$ws1.Range("E21").Value = 213.242846720866
$ws1.Range("F21").Value = [double]"4.9534089201199799E-47"
$ws1.Range("G21").Value = [double]"1.86330224455971E-47"
$ws1.Range("H21").Value = 10.9066775032879

# Row 22 ("Null model", K=3) - values only
$ws1.Range("D22").Value = -7.4448560263490497
$ws1.Range("E22").Value = 219.33548363212699
$ws1.Range("F22").Value = [double]"2.3545335682256701E-48"
$ws1.Range("G22").Value = [double]"8.8569463036776904E-49"
$ws1.Range("H22").Value = 6.8046197939964399

# --- Window / selection state: AIC tab becomes the active tab/sheet,
#     with B18:I22 selected (instead of modelParameters / K11:R12) ---
$ws1.Select() | Out-Null
$ws1.Range("B18:I22").Select() | Out-Null
